# Refresh the cryptocurrency Price (column D) and Volume(1h) change
# (column E) figures to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.831.03'
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").Value = '3.525.74'
$ws.Range("E3").Value = '  +0.93%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '196.78'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.631'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.50%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.198'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.648'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.70'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000301'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.51'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("D14").Value = '4.084.79'
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '596.40'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.72%  '
$ws.Range("D16").Value = '70.032.43'
$ws.Range("E16").Value = '  +0.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.72'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.96%  '
$ws.Range("D19").Value = '3.533.95'
$ws.Range("E19").Value = '  +1.14%  '
$ws.Range("E20").Value = '  +1.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.992'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '101.53'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.42%  '
$ws.Range("E25").Value = '  -1.17%  '
$ws.Range("E26").Value = '  +5.88%  '
$ws.Range("E27").Value = '  -0.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.36'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +11.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.08'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.45'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("E33").Value = '  -1.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.10'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.43%  '
$ws.Range("D35").Value = '0.0₃0858'
$ws.Range("E35").Value = '  +10.08%  '
$ws.Range("D36").Value = '3.705.93'
$ws.Range("E36").Value = '  +2.56%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").Value = '  -4.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.62'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.392'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '490.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.132'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0453'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.37%  '
$ws.Range("E46").Value = '  -2.36%  '
$ws.Range("E47").Value = '  -0.65%  '
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.57'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000251'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.83'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.22%  '

Write-Host "Applied crypto price/volume updates"
